$d = $word.ActiveDocument

$d.Content.Find.Execute("23+14=", $true, $false, $false, $false, $false, $true, 1, $false, "61-45=", 2) | Out-Null
$d.Content.Find.Execute("39+5=", $true, $false, $false, $false, $false, $true, 1, $false, "29-22=", 2) | Out-Null
$d.Content.Find.Execute("88-46=", $true, $false, $false, $false, $false, $true, 1, $false, "11-9=", 2) | Out-Null
$d.Content.Find.Execute("54-33=", $true, $false, $false, $false, $false, $true, 1, $false, "70-36=", 2) | Out-Null
$d.Content.Find.Execute("71-44=", $true, $false, $false, $false, $false, $true, 1, $false, "49+5=", 2) | Out-Null
$d.Content.Find.Execute("5+7=", $true, $false, $false, $false, $false, $true, 1, $false, "96-18=", 2) | Out-Null
$d.Content.Find.Execute("79+12=", $true, $false, $false, $false, $false, $true, 1, $false, "20+41=", 2) | Out-Null
$d.Content.Find.Execute("70-17=", $true, $false, $false, $false, $false, $true, 1, $false, "85-66=", 2) | Out-Null
$d.Content.Find.Execute("13+44=", $true, $false, $false, $false, $false, $true, 1, $false, "44+41=", 2) | Out-Null
$d.Content.Find.Execute("47-26=", $true, $false, $false, $false, $false, $true, 1, $false, "77+11=", 2) | Out-Null
$d.Content.Find.Execute("80-73=", $true, $false, $false, $false, $false, $true, 1, $false, "69+9=", 2) | Out-Null
$d.Content.Find.Execute("51-33=", $true, $false, $false, $false, $false, $true, 1, $false, "32+35=", 2) | Out-Null
$d.Content.Find.Execute("7+80=", $true, $false, $false, $false, $false, $true, 1, $false, "38+36=", 2) | Out-Null
$d.Content.Find.Execute("61+24=", $true, $false, $false, $false, $false, $true, 1, $false, "85-76=", 2) | Out-Null
$d.Content.Find.Execute("40-26=", $true, $false, $false, $false, $false, $true, 1, $false, "54+17=", 2) | Out-Null
$d.Content.Find.Execute("44-11=", $true, $false, $false, $false, $false, $true, 1, $false, "58+1=", 2) | Out-Null
$d.Content.Find.Execute("30+61=", $true, $false, $false, $false, $false, $true, 1, $false, "51+22=", 2) | Out-Null
$d.Content.Find.Execute("36+19=", $true, $false, $false, $false, $false, $true, 1, $false, "30-21=", 2) | Out-Null
$d.Content.Find.Execute("27+70=", $true, $false, $false, $false, $false, $true, 1, $false, "22+5=", 2) | Out-Null
$d.Content.Find.Execute("79-44=", $true, $false, $false, $false, $false, $true, 1, $false, "4+19=", 2) | Out-Null
$d.Content.Find.Execute("68+0=", $true, $false, $false, $false, $false, $true, 1, $false, "81-19=", 2) | Out-Null
$d.Content.Find.Execute("95-48=", $true, $false, $false, $false, $false, $true, 1, $false, "97-17=", 2) | Out-Null
$d.Content.Find.Execute("72-12=", $true, $false, $false, $false, $false, $true, 1, $false, "69-61=", 2) | Out-Null
$d.Content.Find.Execute("7+21=", $true, $false, $false, $false, $false, $true, 1, $false, "7+30=", 2) | Out-Null
$d.Content.Find.Execute("36-34=", $true, $false, $false, $false, $false, $true, 1, $false, "27+13=", 2) | Out-Null
$d.Content.Find.Execute("59-52=", $true, $false, $false, $false, $false, $true, 1, $false, "80-57=", 2) | Out-Null
$d.Content.Find.Execute("48-46=", $true, $false, $false, $false, $false, $true, 1, $false, "29-19=", 2) | Out-Null
$d.Content.Find.Execute("4+36=", $true, $false, $false, $false, $false, $true, 1, $false, "12+80=", 2) | Out-Null
$d.Content.Find.Execute("44-43=", $true, $false, $false, $false, $false, $true, 1, $false, "72-6=", 2) | Out-Null
$d.Content.Find.Execute("46+44=", $true, $false, $false, $false, $false, $true, 1, $false, "97-11=", 2) | Out-Null
$d.Content.Find.Execute("26-7=", $true, $false, $false, $false, $false, $true, 1, $false, "39+2=", 2) | Out-Null
$d.Content.Find.Execute("70-68=", $true, $false, $false, $false, $false, $true, 1, $false, "30+62=", 2) | Out-Null
$d.Content.Find.Execute("62+14=", $true, $false, $false, $false, $false, $true, 1, $false, "4+59=", 2) | Out-Null
$d.Content.Find.Execute("80-69=", $true, $false, $false, $false, $false, $true, 1, $false, "89-87=", 2) | Out-Null
$d.Content.Find.Execute("72+19=", $true, $false, $false, $false, $false, $true, 1, $false, "3+81=", 2) | Out-Null
$d.Content.Find.Execute("92-75=", $true, $false, $false, $false, $false, $true, 1, $false, "25+33=", 2) | Out-Null
$d.Content.Find.Execute("11+38=", $true, $false, $false, $false, $false, $true, 1, $false, "17-10=", 2) | Out-Null
$d.Content.Find.Execute("48-2=", $true, $false, $false, $false, $false, $true, 1, $false, "98+1=", 2) | Out-Null
$d.Content.Find.Execute("38+52=", $true, $false, $false, $false, $false, $true, 1, $false, "13+39=", 2) | Out-Null
$d.Content.Find.Execute("31+32=", $true, $false, $false, $false, $false, $true, 1, $false, "14+10=", 2) | Out-Null
$d.Content.Find.Execute("10+2=", $true, $false, $false, $false, $false, $true, 1, $false, "19+80=", 2) | Out-Null
$d.Content.Find.Execute("98-90=", $true, $false, $false, $false, $false, $true, 1, $false, "44+5=", 2) | Out-Null
$d.Content.Find.Execute("21+40=", $true, $false, $false, $false, $false, $true, 1, $false, "19+27=", 2) | Out-Null
$d.Content.Find.Execute("37+5=", $true, $false, $false, $false, $false, $true, 1, $false, "18+67=", 2) | Out-Null
$d.Content.Find.Execute("25-0=", $true, $false, $false, $false, $false, $true, 1, $false, "9+0=", 2) | Out-Null
$d.Content.Find.Execute("19+43=", $true, $false, $false, $false, $false, $true, 1, $false, "82-27=", 2) | Out-Null
$d.Content.Find.Execute("89-60=", $true, $false, $false, $false, $false, $true, 1, $false, "4+32=", 2) | Out-Null
$d.Content.Find.Execute("75+14=", $true, $false, $false, $false, $false, $true, 1, $false, "26+8=", 2) | Out-Null
$d.Content.Find.Execute("93-9=", $true, $false, $false, $false, $false, $true, 1, $false, "76-3=", 2) | Out-Null
$d.Content.Find.Execute("48-1=", $true, $false, $false, $false, $false, $true, 1, $false, "91-6=", 2) | Out-Null
$d.Content.Find.Execute("52+33=", $true, $false, $false, $false, $false, $true, 1, $false, "82-14=", 2) | Out-Null
$d.Content.Find.Execute("59+0=", $true, $false, $false, $false, $false, $true, 1, $false, "32+50=", 2) | Out-Null
$d.Content.Find.Execute("46+49=", $true, $false, $false, $false, $false, $true, 1, $false, "82+6=", 2) | Out-Null
$d.Content.Find.Execute("16+77=", $true, $false, $false, $false, $false, $true, 1, $false, "39+4=", 2) | Out-Null
$d.Content.Find.Execute("86-58=", $true, $false, $false, $false, $false, $true, 1, $false, "70-58=", 2) | Out-Null
$d.Content.Find.Execute("48-43=", $true, $false, $false, $false, $false, $true, 1, $false, "3+6=", 2) | Out-Null
$d.Content.Find.Execute("85-8=", $true, $false, $false, $false, $false, $true, 1, $false, "29-17=", 2) | Out-Null
$d.Content.Find.Execute("66-20=", $true, $false, $false, $false, $false, $true, 1, $false, "55+20=", 2) | Out-Null
$d.Content.Find.Execute("37-28=", $true, $false, $false, $false, $false, $true, 1, $false, "37-20=", 2) | Out-Null
$d.Content.Find.Execute("10+60=", $true, $false, $false, $false, $false, $true, 1, $false, "10+39=", 2) | Out-Null
$d.Content.Find.Execute("70-60=", $true, $false, $false, $false, $false, $true, 1, $false, "6+44=", 2) | Out-Null
$d.Content.Find.Execute("18+72=", $true, $false, $false, $false, $false, $true, 1, $false, "28-22=", 2) | Out-Null
$d.Content.Find.Execute("10-6=", $true, $false, $false, $false, $false, $true, 1, $false, "45-11=", 2) | Out-Null
$d.Content.Find.Execute("36-28=", $true, $false, $false, $false, $false, $true, 1, $false, "49-38=", 2) | Out-Null
$d.Content.Find.Execute("71-23=", $true, $false, $false, $false, $false, $true, 1, $false, "61-29=", 2) | Out-Null
$d.Content.Find.Execute("56+1=", $true, $false, $false, $false, $false, $true, 1, $false, "9+38=", 2) | Out-Null
$d.Content.Find.Execute("43+52=", $true, $false, $false, $false, $false, $true, 1, $false, "80+9=", 2) | Out-Null
$d.Content.Find.Execute("6+33=", $true, $false, $false, $false, $false, $true, 1, $false, "19-14=", 2) | Out-Null
$d.Content.Find.Execute("14+24=", $true, $false, $false, $false, $false, $true, 1, $false, "60+31=", 2) | Out-Null
$d.Content.Find.Execute("93-84=", $true, $false, $false, $false, $false, $true, 1, $false, "87+1=", 2) | Out-Null
$d.Content.Find.Execute("11+74=", $true, $false, $false, $false, $false, $true, 1, $false, "84-38=", 2) | Out-Null
$d.Content.Find.Execute("31+10=", $true, $false, $false, $false, $false, $true, 1, $false, "6+70=", 2) | Out-Null
$d.Content.Find.Execute("47-22=", $true, $false, $false, $false, $false, $true, 1, $false, "72-32=", 2) | Out-Null
$d.Content.Find.Execute("97-23=", $true, $false, $false, $false, $false, $true, 1, $false, "18+64=", 2) | Out-Null
$d.Content.Find.Execute("22-1=", $true, $false, $false, $false, $false, $true, 1, $false, "86-15=", 2) | Out-Null
$d.Content.Find.Execute("60-33=", $true, $false, $false, $false, $false, $true, 1, $false, "65+7=", 2) | Out-Null
$d.Content.Find.Execute("83-56=", $true, $false, $false, $false, $false, $true, 1, $false, "65-35=", 2) | Out-Null
$d.Content.Find.Execute("87-18=", $true, $false, $false, $false, $false, $true, 1, $false, "74-66=", 2) | Out-Null
$d.Content.Find.Execute("50-7=", $true, $false, $false, $false, $false, $true, 1, $false, "40-32=", 2) | Out-Null
$d.Content.Find.Execute("15-9=", $true, $false, $false, $false, $false, $true, 1, $false, "47+14=", 2) | Out-Null
$d.Content.Find.Execute("28-26=", $true, $false, $false, $false, $false, $true, 1, $false, "1+6=", 2) | Out-Null
$d.Content.Find.Execute("28+60=", $true, $false, $false, $false, $false, $true, 1, $false, "16+55=", 2) | Out-Null
$d.Content.Find.Execute("69-44=", $true, $false, $false, $false, $false, $true, 1, $false, "73-31=", 2) | Out-Null
$d.Content.Find.Execute("7+7=", $true, $false, $false, $false, $false, $true, 1, $false, "33-4=", 2) | Out-Null
$d.Content.Find.Execute("91-29=", $true, $false, $false, $false, $false, $true, 1, $false, "47-27=", 2) | Out-Null
$d.Content.Find.Execute("35+34=", $true, $false, $false, $false, $false, $true, 1, $false, "4+29=", 2) | Out-Null
$d.Content.Find.Execute("75-17=", $true, $false, $false, $false, $false, $true, 1, $false, "84-43=", 2) | Out-Null
$d.Content.Find.Execute("85+11=", $true, $false, $false, $false, $false, $true, 1, $false, "53+7=", 2) | Out-Null
$d.Content.Find.Execute("64-21=", $true, $false, $false, $false, $false, $true, 1, $false, "91+3=", 2) | Out-Null
$d.Content.Find.Execute("90-28=", $true, $false, $false, $false, $false, $true, 1, $false, "82-41=", 2) | Out-Null
$d.Content.Find.Execute("23+10=", $true, $false, $false, $false, $false, $true, 1, $false, "69+12=", 2) | Out-Null
$d.Content.Find.Execute("42+1=", $true, $false, $false, $false, $false, $true, 1, $false, "76-57=", 2) | Out-Null
$d.Content.Find.Execute("59-44=", $true, $false, $false, $false, $false, $true, 1, $false, "52-26=", 2) | Out-Null
$d.Content.Find.Execute("29-15=", $true, $false, $false, $false, $false, $true, 1, $false, "92-34=", 2) | Out-Null
$d.Content.Find.Execute("89-70=", $true, $false, $false, $false, $false, $true, 1, $false, "45-6=", 2) | Out-Null
$d.Content.Find.Execute("75-38=", $true, $false, $false, $false, $false, $true, 1, $false, "54-34=", 2) | Out-Null
$d.Content.Find.Execute("36-10=", $true, $false, $false, $false, $false, $true, 1, $false, "82-40=", 2) | Out-Null
$d.Content.Find.Execute("3+12=", $true, $false, $false, $false, $false, $true, 1, $false, "27-1=", 2) | Out-Null
$d.Content.Find.Execute("78-2=", $true, $false, $false, $false, $false, $true, 1, $false, "92+5=", 2) | Out-Null
$d.Content.Find.Execute("66-36=", $true, $false, $false, $false, $false, $true, 1, $false, "73-35=", 2) | Out-Null
